$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order used for each data row: B, C, D, E, F, I, J, K, L, M
$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13)

$data = @(
    @(1.02, 1.072518403766447, 1.071734474037007, 1.076252960779337, 1.086234518589867, 1.061579896306071, 1.077437545509292, 1.074430858838986, 1.078937380886507, 1.088892888960648),
    @(1.02, 1.073887812894446, 1.072835757304793, 1.077505242709078, 1.087588028979667, 1.062110832322893, 1.078463321328205, 1.075348049449879, 1.080006068472969, 1.090064456341429),
    @(1.02, 1.0747725579955, 1.073547145343368, 1.078314015975195, 1.088462685324939, 1.062452409196289, 1.07912521098819, 1.075939723953691, 1.080695480968981, 1.090820836773028),
    @(1.02, 1.07514418642054, 1.073845925778152, 1.078653661465729, 1.088830119571591, 1.062595538018972, 1.079403029489041, 1.076188034622839, 1.080984812066548, 1.091138416554309),
    @(1.02, 1.075206565963696, 1.073896075607803, 1.078710668312866, 1.088891797644499, 1.062619542485075, 1.079449650768315, 1.076229702013602, 1.081033362971674, 1.091191716133434),
    @(1.02, 1.074777524956682, 1.073551138785309, 1.07831855575399, 1.088467596058662, 1.062454323534629, 1.079128924935704, 1.075943043575093, 1.080699348975789, 1.090825081865265),
    @(1.02, 1.072981485868488, 1.072106912035172, 1.076676496643899, 1.086692186202848, 1.061759738929663, 1.077784598173005, 1.074741204604064, 1.079298986533322, 1.089289181016534),
    @(1.02, 1.069806004645032, 1.069552508372834, 1.073770976952369, 1.083554599238295, 1.060520557561418, 1.075401311360128, 1.072609382870747, 1.076815071582797, 1.086569471428185),
    @(1.02, 1.067681504702729, 1.067842937746212, 1.07182556937016, 1.081456427884369, 1.059684032647646, 1.073802486140463, 1.071178494396457, 1.075147868078264, 1.084747108415521),
    @(1.02, 1.066759709130937, 1.067101044037986, 1.070981122375452, 1.080546291927054, 1.059319304154323, 1.073107750589969, 1.070556554809426, 1.074423211847765, 1.083955747728346),
    @(1.02, 1.066417024672261, 1.066825220532878, 1.07066713938956, 1.080207977841416, 1.059183447806929, 1.072849323771709, 1.070325180092355, 1.074153623743957, 1.083661454673405),
    @(1.02, 1.066490544823996, 1.066884397056838, 1.070734504288632, 1.08028055873841, 1.059212606699472, 1.072904774082825, 1.070374827062023, 1.074211470341856, 1.08372459730911),
    @(1.02, 1.066731388654239, 1.067078249533466, 1.070955174963404, 1.080518331898399, 1.059308081995592, 1.073086396550685, 1.070537436647544, 1.074400936183085, 1.083931428475725),
    @(1.02, 1.066879742038352, 1.067197655012169, 1.071091095149669, 1.080664798587899, 1.059366857034231, 1.073198250808612, 1.070637578109976, 1.074517616705373, 1.084058817891281),
    @(1.02, 1.067742641167737, 1.067892139843757, 1.071881568323865, 1.081516796045177, 1.059708185370367, 1.073848541710123, 1.071219720444575, 1.075195902761956, 1.084799580167482),
    @(1.02, 1.068283408497727, 1.068327329676091, 1.072376851859117, 1.082050794807497, 1.05992161799101, 1.074255796597578, 1.071584248966395, 1.075620634289882, 1.085263629947494),
    @(1.02, 1.068598648692316, 1.06858101080674, 1.072665542816396, 1.082362111980111, 1.060045867889657, 1.074493107040403, 1.071796645346101, 1.075868108448337, 1.085534084162108),
    @(1.02, 1.068706107095045, 1.068667482912435, 1.072763945317339, 1.082468236952685, 1.060088192987353, 1.074573984167789, 1.071869028673829, 1.075952446046706, 1.085626265238114),
    @(1.02, 1.068225407977326, 1.068280654278548, 1.072323733311102, 1.081993517903728, 1.059898743718733, 1.074212126292803, 1.071545162015132, 1.075575092025227, 1.085213864409867),
    @(1.02, 1.066660474171746, 1.067021171784082, 1.070890201773278, 1.080448320557787, 1.059279977414661, 1.073032923549055, 1.070489562116989, 1.074345154816122, 1.083870531465472),
    @(1.02, 1.065674865204818, 1.066227829573951, 1.069987041873091, 1.079475349333228, 1.058888734635827, 1.072289362166281, 1.06982378705905, 1.073569420110346, 1.083023917164291),
    @(1.02, 1.066197516137465, 1.066648534942497, 1.070466000738367, 1.079991278906504, 1.059096349427043, 1.072683743736635, 1.070176925543545, 1.073980883550572, 1.083472915763463),
    @(1.02, 1.068251616472297, 1.068301745368641, 1.072347735912629, 1.082019399350404, 1.059909080364469, 1.074231859747773, 1.071562824426811, 1.075595671431296, 1.08523635198788),
    @(1.02, 1.070628239895182, 1.070214033709291, 1.074523579055632, 1.084366852156093, 1.060842737544397, 1.07601918444925, 1.073162197121179, 1.077459185516621, 1.087274184320543)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $rowValues[$j]
    }
}

Write-Output "Updated $($data.Count) rows x $($cols.Count) cols"
